# BM_Tru.NghiemThu.docx edit:
# The "Đơn vị thi công" (construction-unit) table — the table whose rows
# carry the {#manager}/{#leaders} merge fields — had its overall width
# switched from "auto" to a fixed 8651 dxa, and its third (rightmost)
# column widened from 3645 dxa to 3933 dxa (so both table cells in that
# column grow to match).
#
# Word's Table/Column PreferredWidth API takes its value in points and
# always persists it as twentieths-of-a-point (dxa) in the XML, so to
# land on an exact dxa figure we divide the target dxa value by 20
# before assigning it.

$d = $word.ActiveDocument

# Locate the correct table: the one containing the {#manager} merge
# field (the "Đơn vị thi công" block), rather than relying on a bare
# positional index.
$target = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Range.Text -like "*{#manager}*") {
        $target = $candidate
        break
    }
}

# Set the overall table width: <w:tblW w:w="8651" w:type="dxa"/>
$target.PreferredWidthType = 3
$target.PreferredWidth = 8651 / 20

# Set the third column's width: <w:gridCol w:w="3933"/> /
# <w:tcW w:w="3933" w:type="dxa"/> on every cell in that column.
$col = $target.Columns.Item(3)
$col.PreferredWidthType = 3
$col.PreferredWidth = 3933 / 20
